$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column C
$ws.Range("C1").Value = "GroepUUID"

# Replace groepN values with UUIDs
$ws.Range("C2").Value = "972fa507-8021-4933-9646-3e9c238f83ec"
$ws.Range("C3").Value = "bcda9a3f-ab3f-4135-ba66-ecbd03fb745c"
$ws.Range("C4").Value = "bcda9a3f-ab3f-4135-ba66-ecbd03fb745c"
$ws.Range("C5").Value = "88d796c9-d588-4b88-956e-73ad02eb5ea7"
$ws.Range("C6").Value = "a1308055-c284-483b-9a6a-50b357cbbcd1"
$ws.Range("C7").Value = "e5a57688-313e-4a24-bf31-fcace58880a0"
$ws.Range("C8").Value = "5a5e9758-2bac-4f05-8f18-bffdbbff5ca1"
$ws.Range("C9").Value = "57b61f85-cc13-4b47-bd55-ad4990f2818f"

# Column width adjustment (bestFit width becomes 38 due to new longer content).
# ColumnWidth round-trips through the engine's internal storage with a
# constant +0.8333... padding offset, so back the request off by that
# amount to land exactly on the target stored width of 38.
$ws.Columns("C").ColumnWidth = 37.1666666666667

# Update selection to B2
$ws.Range("B2").Select()
